$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 (dates + volume/price columns)
$ws.Range("D2").Value = 44425
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 13000
$ws.Range("P2").Value = 1300

$ws.Range("D3").Value = 44348
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 1000
